$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("L2").Value = "GPT-0532"
$ws.Range("M2").Value = "OKI TONER M C3100/3000/3200/5100/5150/5200/5300/5400/5510 MAGENTA ΣΥΜΒΑΤΟ 3000 ΣΕΛΙΔΕΣ"
$ws.Range("N2").Value = 7
$ws.Range("O2").Value = "'4€"
$ws.Range("O2").Style = "Normal"

# Remove row 3 entirely (shifts dimension back to A1:O2)
$ws.Rows("3:3").Delete()
